$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Update "总计" (summary) sheet: insert a new row for 2022-Q3 at the top
#    of the data (row 2), shifting all existing rows down by one.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Rows.Item(2).Insert()
# Copy the (soon to be old-row-2, now row-3) formatting pattern onto the
# newly inserted row so the column-A style (bordered) and plain B:D cells
# match the rest of the table.
$summary.Range("A3:D3").Copy()
$summary.Range("A2:D2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 12
$summary.Range("D2").Value = 0.32

# ---------------------------------------------------------------------------
# 2) Insert a brand-new "2022-Q3" sheet before the existing "2022-Q2" sheet,
#    built from a copy of "2022-Q2" so it inherits identical formatting.
# ---------------------------------------------------------------------------
$q2 = $wb.Worksheets.Item(2)
$q2.Copy($q2)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# Wipe the copied data rows (2-7), then grow the table to 12 data rows
# (rows 2-13) by duplicating the last row's formatting downward.
$q3.Range("A2:H7").ClearContents()
$q3.Range("A7:H7").Copy()
$q3.Range("A8:H8").PasteSpecial(-4122)
$q3.Range("A9:H9").PasteSpecial(-4122)
$q3.Range("A10:H10").PasteSpecial(-4122)
$q3.Range("A11:H11").PasteSpecial(-4122)
$q3.Range("A12:H12").PasteSpecial(-4122)
$q3.Range("A13:H13").PasteSpecial(-4122)

$data = @(
    @(0,  "001277", "博时国企改革主题股票A",       "2.42", "85.79", "3.39", "0.0820", 5),
    @(1,  "001541", "汇添富民营新动力股票",         "2.21", "91.50", "3.39", "0.0749", 7),
    @(2,  "050014", "博时创业成长混合A",           "1.47", "82.43", "3.72", "0.0547", 7),
    @(3,  "011269", "中银证券优势制造股票A",        "0.76", "92.92", "4.62", "0.0351", 9),
    @(4,  "011270", "中银证券优势制造股票C",        "0.58", "92.92", "4.62", "0.0268", 9),
    @(5,  "160519", "博时睿利事件驱动灵活配置混合",  "0.55", "82.58", "4.03", "0.0222", 6),
    @(6,  "002149", "嘉实新优选灵活配置混合",        "0.18", "91.79", "4.90", "0.0088", 8),
    @(7,  "620004", "金元顺安价值增长混合",         "0.34", "74.99", "1.92", "0.0065", 7),
    @(8,  "004913", "中银证券聚瑞混合A",           "0.07", "43.92", "2.88", "0.0020", 8),
    @(9,  "002553", "博时创业成长混合C",           "0.05", "82.43", "3.72", "0.0019", 7),
    @(10, "004914", "中银证券聚瑞混合C",           "0.05", "43.92", "2.88", "0.0014", 8),
    @(11, "014382", "博时国企改革主题股票C",        "0.00", "85.79", "3.39", "0",      5)
)

$r = 2
foreach ($row in $data) {
    $q3.Range("A$r").Value = $row[0]
    $q3.Range("B$r").Value = "'" + $row[1]
    $q3.Range("C$r").Value = $row[2]
    $q3.Range("D$r").Value = "'" + $row[3]
    $q3.Range("E$r").Value = "'" + $row[4]
    $q3.Range("F$r").Value = "'" + $row[5]
    $q3.Range("G$r").Value = "'" + $row[6]
    $q3.Range("H$r").Value = $row[7]
    $q3.Range(("B" + $r + ":G" + $r)).ClearFormats()
    $r++
}

# Row 13's "持有市值(亿元)" (G13) is a genuine numeric 0, not text, in the
# target data - fix it up after the text-forcing loop above.
$q3.Range("G13").Value = 0
